$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "60.266.49"
$c.ClearFormats()
$ws.Range("E2").Value = "  -4.30%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.903.42"
$c.ClearFormats()
$ws.Range("E3").Value = "  -3.69%  "
$ws.Range("E4").Value = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "527.15"
$c.ClearFormats()
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "141.86"
$c.ClearFormats()
$ws.Range("E6").Value = "  -7.77%  "
$ws.Range("E7").Value = "  +0.07%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.549"
$c.ClearFormats()
$ws.Range("E8").Value = "  -2.73%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.904.40"
$c.ClearFormats()
$ws.Range("E9").Value = "  -4.03%  "
$ws.Range("E10").Value = "  -5.35%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.85"
$c.ClearFormats()
$ws.Range("E11").Value = "  -8.80%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.352"
$c.ClearFormats()
$ws.Range("E12").Value = "  -3.53%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "3.412.68"
$c.ClearFormats()
$ws.Range("E13").Value = "  -3.98%  "
$ws.Range("E14").Value = "  +1.18%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "60.435.25"
$c.ClearFormats()
$ws.Range("E15").Value = "  -4.23%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "22.59"
$c.ClearFormats()
$ws.Range("E16").Value = "  -5.70%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.908.84"
$c.ClearFormats()
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("E18").Value = "  -6.40%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.91"
$c.ClearFormats()
$ws.Range("E19").Value = "  -4.00%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.53"
$c.ClearFormats()
$ws.Range("E20").Value = "  -2.90%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "358.92"
$c.ClearFormats()
$ws.Range("E21").Value = "  -9.34%  "
$ws.Range("E22").Value = "  -0.72%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -1.88%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "63.21"
$c.ClearFormats()
$ws.Range("E25").Value = "  -3.11%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.030.37"
$c.ClearFormats()
$ws.Range("E26").Value = "  -3.85%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.449"
$c.ClearFormats()
$ws.Range("E27").Value = "  -3.71%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.176"
$c.ClearFormats()
$ws.Range("E28").Value = "  -6.57%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E29").Value = "  +0.65%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0₃0854"
$c.ClearFormats()
$ws.Range("E30").Value = "  -12.47%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.ClearFormats()
$ws.Range("E31").Value = "  -11.62%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  -4.99%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "19.63"
$c.ClearFormats()
$ws.Range("E34").Value = "  -4.13%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "152.73"
$c.ClearFormats()
$ws.Range("E35").Value = "  -4.51%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.30"
$c.ClearFormats()
$ws.Range("E36").Value = "  -8.26%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.55"
$c.ClearFormats()
$ws.Range("E37").Value = "  -8.20%  "
$ws.Range("E38").Value = "  -9.46%  "
$ws.Range("E39").Value = "  -8.02%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "37.78"
$c.ClearFormats()
$ws.Range("E40").Value = "  +0.42%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.330.50"
$c.ClearFormats()
$ws.Range("E41").Value = "  -7.39%  "
$ws.Range("E42").Value = "  -7.95%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.66"
$c.ClearFormats()
$ws.Range("E43").Value = "  -6.70%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.642"
$c.ClearFormats()
$ws.Range("E44").Value = "  -3.90%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "20.80"
$c.ClearFormats()
$ws.Range("E45").Value = "  -8.12%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0567"
$c.ClearFormats()
$ws.Range("E46").Value = "  -5.26%  "
$ws.Range("E47").Value = "  -0.13%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "4.81"
$c.ClearFormats()
$ws.Range("E48").Value = "  -5.41%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "10.37"
$c.ClearFormats()
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("E50").Value = "  -6.04%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0922"
$c.ClearFormats()
$ws.Range("E51").Value = "  -2.65%  "
